# "scale-parametre indført. kan ikke køre lige pt"
#
# - Introduce a scale parameter: reduce the "endofpipe" sheet's G-column
#   values (G2: 10 -> 2, G3: 6.5 -> 6).
# - The user had last navigated to cell F26 on "inputdisp" (no longer the
#   active tab) and cell G4 on "endofpipe", which becomes the active tab.

$wb = $excel.ActiveWorkbook

# Leave "inputdisp" with the cursor on F26, and make it the inactive tab.
$wsInput = $wb.Worksheets.Item("inputdisp")
$wsInput.Activate()
$wsInput.Range("F26").Select()

# Update the scale parameter values on "endofpipe" and make it the active tab,
# with the cursor left on G4.
$wsEndOfPipe = $wb.Worksheets.Item("endofpipe")
$wsEndOfPipe.Activate()
$wsEndOfPipe.Range("G2").Value = 2
$wsEndOfPipe.Range("G3").Value = 6
$wsEndOfPipe.Range("G4").Select()
